# Adds the new "Cartografia-social-api-front" (frontend) rows to the
# analysis sheet, per commit "[#95] - add frontend datas".
#
# The sheet already holds rows for 5 repositories (rows 2-31). This appends
# 14 more rows (32-45) for the new "Cartografia-social-api-front" repo,
# covering versions MVP, v0.0.1-v0.0.9 (minus v0.0.7), v1.0.0, v1.0.1 and
# the new v1.1.0-v1.1.2 releases. Columns are:
#   A..F = m1..m6 metrics, G = repository, H = version, I = ncloc,
#   J = code_quality, K = testing_status, L = Maintainability,
#   M = Reliability, N = total

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$repo = "Cartografia-social-api-front"

$newRows = @(
    @(0.8181818181818182,  0,                    1,                    1,                    1,                    1,                    $repo, "MVP",    1921, 0.6000000000000001, 1,                    0.3,                  0.5,                  0.8),
    @(0.6909090909090909,  0.03636363636363636,  1,                    1,                    1,                    1,                    $repo, "v0.0.1", 1412, 0.5700000000000001, 1,                    0.285,                0.5,                  0.785),
    @(0.6909090909090909,  0.03636363636363636,  1,                    1,                    1,                    1,                    $repo, "v0.0.2", 1412, 0.5700000000000001, 1,                    0.285,                0.5,                  0.785),
    @(0.6,                 0,                    1,                    1,                    1,                    1,                    $repo, "v0.0.3", 525,  0.528,              1,                    0.264,                0.5,                  0.764),
    @(0.6,                 0,                    1,                    1,                    1,                    1,                    $repo, "v0.0.4", 525,  0.528,              1,                    0.264,                0.5,                  0.764),
    @(0.7222222222222222,  0,                    0.9444444444444444,   0.9444444444444444,   0.9444444444444444,   0.9444444444444444,   $repo, "v0.0.5", 975,  0.55,               0.9444444444444444,   0.275,                0.4722222222222222,   0.7472222222222222),
    @(0.75,                0,                    0.95,                 0.95,                 0.95,                 0.95,                 $repo, "v0.0.6", 1077, 0.5609999999999999, 0.95,                 0.2805,               0.475,                0.7554999999999999),
    @(0.8148148148148148,  0,                    1,                    1,                    1,                    1,                    $repo, "v0.0.8", 1625, 0.5988888888888889, 1,                    0.2994444444444445,  0.5,                  0.7994444444444444),
    @(0.8333333333333334,  0,                    1,                    1,                    1,                    1,                    $repo, "v0.0.9", 1788, 0.605,              1,                    0.3025,               0.5,                  0.8025),
    @(0.8285714285714286,  0,                    0.9428571428571428,   0.9428571428571428,   0.9428571428571428,   0.9428571428571428,   $repo, "v1.0.0", 2086, 0.5845714285714286, 0.9428571428571428,  0.2922857142857143,  0.4714285714285714,  0.7637142857142858),
    @(0.8333333333333334,  0,                    0.9444444444444444,   0.9444444444444444,   0.9444444444444444,   0.9444444444444444,   $repo, "v1.0.1", 2136, 0.5866666666666667, 0.9444444444444444,  0.2933333333333333,  0.4722222222222222,  0.7655555555555555),
    @(0.7692307692307693,  0,                    0.9230769230769231,   0.9230769230769231,   0.9230769230769231,   0.9230769230769231,   $repo, "v1.1.0", 3091, 0.5584615384615386, 0.9230769230769231,  0.2792307692307693,  0.4615384615384616,  0.7407692307692308),
    @(0.7857142857142857,  0,                    0.9285714285714286,   0.9285714285714286,   0.9285714285714286,   0.9285714285714286,   $repo, "v1.1.1", 3545, 0.5657142857142857, 0.9285714285714286,  0.2828571428571429,  0.4642857142857143,  0.7471428571428571),
    @(0.7857142857142857,  0,                    0.9285714285714286,   0.9285714285714286,   0.9285714285714286,   0.9285714285714286,   $repo, "v1.1.2", 3511, 0.5657142857142857, 0.9285714285714286,  0.2828571428571429,  0.4642857142857143,  0.7471428571428571)
)

$firstNewRow = 32
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $excelRow = $firstNewRow + $i
    $rowValues = $newRows[$i]
    for ($col = 0; $col -lt $rowValues.Length; $col++) {
        $ws.Cells.Item($excelRow, $col + 1).Value = $rowValues[$col]
    }
}
